# Adding generated reports from pipeline run
#
# All_Components: the Flask framework row (68) loses its License/Hashes/
# Description values, and 5 new Flask-transitive Python packages are
# inserted right after it (pushing the existing GitHub Action rows from
# 69-72 down to 74-77).
#
# Dependency_Mapping: 9 new direct/transitive dependency edges for the
# Flask app are appended after the existing last row (68).

$wb = $excel.ActiveWorkbook

$components = $wb.Worksheets.Item("All_Components")
$mapping = $wb.Worksheets.Item("Dependency_Mapping")

# ---------------------------------------------------------------------
# All_Components sheet
# ---------------------------------------------------------------------

# Row 68 (flask) keeps Group/Name/Version/PURL/Type, but License, Hashes
# and Description are cleared out.
$components.Range("F68").Value = ""
$components.Range("G68").Value = ""
$components.Range("H68").Value = ""

# New rows 69-73: Flask's transitive Python dependencies.
$newComponentRows = @(
    @{ Row = 69; Group = ""; Name = "jinja2";       Version = "3.1.6"; Purl = "pkg:pypi/jinja2@3.1.6";       Type = "library" },
    @{ Row = 70; Group = ""; Name = "markupsafe";   Version = "3.0.2"; Purl = "pkg:pypi/markupsafe@3.0.2";   Type = "library" },
    @{ Row = 71; Group = ""; Name = "werkzeug";     Version = "3.1.3"; Purl = "pkg:pypi/werkzeug@3.1.3";     Type = "library" },
    @{ Row = 72; Group = ""; Name = "click";        Version = "8.1.8"; Purl = "pkg:pypi/click@8.1.8";        Type = "library" },
    @{ Row = 73; Group = ""; Name = "itsdangerous"; Version = "2.2.0"; Purl = "pkg:pypi/itsdangerous@2.2.0"; Type = "library" }
)

foreach ($row in $newComponentRows) {
    $components.Range("A" + $row.Row).Value = $row.Group
    $components.Range("B" + $row.Row).Value = $row.Name
    $components.Range("C" + $row.Row).Value = $row.Version
    $components.Range("D" + $row.Row).Value = $row.Purl
    $components.Range("E" + $row.Row).Value = $row.Type
    $components.Range("F" + $row.Row).Value = ""
    $components.Range("G" + $row.Row).Value = ""
    $components.Range("H" + $row.Row).Value = ""
}

# Rows 74-77: the 4 GitHub Action rows that used to be 69-72, now shifted
# down by 5 to make room for the new Flask dependency rows above.
$shiftedComponentRows = @(
    @{ Row = 74; Group = "actions";   Name = "checkout";                                 Version = "v4";                                         Purl = "pkg:github/actions/checkout@v4";                                         Type = "application" },
    @{ Row = 75; Group = "";          Name = "snyk%2Factions%2Fiac";                      Version = "14818c4695ecc4045f33c9cee9e795a788711ca4";    Purl = "pkg:github/snyk%2Factions%2Fiac@14818c4695ecc4045f33c9cee9e795a788711ca4"; Type = "library" },
    @{ Row = 76; Group = "";          Name = "github%2Fcodeql-action%2Fupload-sarif";     Version = "v3";                                         Purl = "pkg:github/github%2Fcodeql-action%2Fupload-sarif@v3";                     Type = "library" },
    @{ Row = 77; Group = "oxsecurity"; Name = "ox-security-scan";                         Version = "main";                                       Purl = "pkg:github/oxsecurity/ox-security-scan@main";                             Type = "library" }
)

foreach ($row in $shiftedComponentRows) {
    $components.Range("A" + $row.Row).Value = $row.Group
    $components.Range("B" + $row.Row).Value = $row.Name
    $components.Range("C" + $row.Row).Value = $row.Version
    $components.Range("D" + $row.Row).Value = $row.Purl
    $components.Range("E" + $row.Row).Value = $row.Type
    $components.Range("F" + $row.Row).Value = ""
    $components.Range("G" + $row.Row).Value = ""
    $components.Range("H" + $row.Row).Value = ""
}

# ---------------------------------------------------------------------
# Dependency_Mapping sheet
# ---------------------------------------------------------------------

$newMappingRows = @(
    @{ Row = 69; Direct = "pkg:pypi/app@latest"; Transitive = "flask" },
    @{ Row = 70; Direct = "pkg:pypi/app@latest"; Transitive = "jinja2" },
    @{ Row = 71; Direct = "pkg:pypi/app@latest"; Transitive = "werkzeug" },
    @{ Row = 72; Direct = "jinja2";              Transitive = "markupsafe" },
    @{ Row = 73; Direct = "werkzeug";            Transitive = "markupsafe" },
    @{ Row = 74; Direct = "flask";               Transitive = "click" },
    @{ Row = 75; Direct = "flask";               Transitive = "itsdangerous" },
    @{ Row = 76; Direct = "flask";               Transitive = "jinja2" },
    @{ Row = 77; Direct = "flask";               Transitive = "werkzeug" }
)

foreach ($row in $newMappingRows) {
    $mapping.Range("A" + $row.Row).Value = $row.Direct
    $mapping.Range("B" + $row.Row).Value = $row.Transitive
}
